$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.308.32"
$ws.Range("D3").Value = "'3.673.43"
$ws.Range("E3").Value = "'  -0.36%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("D5").Value = "'683.08"
$ws.Range("E5").Value = "'  -0.11%  "
$ws.Range("D6").Value = "'157.70"
$ws.Range("E6").Value = "'  -2.97%  "
$ws.Range("E7").Value = "'  -0.05%  "
$ws.Range("E8").Value = "'  -1.18%  "
$ws.Range("E9").Value = "'  -1.97%  "
$ws.Range("E10").Value = "'  -3.97%  "
$ws.Range("D11").Value = "'0.435"
$ws.Range("E11").Value = "'  -3.48%  "
$ws.Range("D12").Value = "'0.0000232"
$ws.Range("E12").Value = "'  -2.25%  "
$ws.Range("D13").Value = "'4.293.27"
$ws.Range("E13").Value = "'  -0.36%  "
$ws.Range("D14").Value = "'32.14"
$ws.Range("E14").Value = "'  -4.41%  "
$ws.Range("D15").Value = "'3.666.16"
$ws.Range("E15").Value = "'  -0.76%  "
$ws.Range("D16").Value = "'69.295.35"
$ws.Range("E16").Value = "'  -0.14%  "
$ws.Range("E17").Value = "'  +1.97%  "
$ws.Range("D18").Value = "'15.81"
$ws.Range("E18").Value = "'  -3.17%  "
$ws.Range("D19").Value = "'6.37"
$ws.Range("E19").Value = "'  -4.26%  "
$ws.Range("D20").Value = "'471.69"
$ws.Range("E20").Value = "'  -2.19%  "
$ws.Range("D21").Value = "'9.92"
$ws.Range("E21").Value = "'  +1.09%  "
$ws.Range("D22").Value = "'0.647"
$ws.Range("E22").Value = "'  -3.13%  "
$ws.Range("D23").Value = "'79.93"
$ws.Range("E23").Value = "'  -0.11%  "
$ws.Range("D24").Value = "'3.819.04"
$ws.Range("E24").Value = "'  -0.36%  "
$ws.Range("E25").Value = "'  -0.14%  "
$ws.Range("E26").Value = "'  -5.45%  "
$ws.Range("D27").Value = "'10.90"
$ws.Range("E27").Value = "'  -5.08%  "
$ws.Range("D28").Value = "'9.11"
$ws.Range("E28").Value = "'  -5.08%  "
$ws.Range("D29").Value = "'2.70"
$ws.Range("E29").Value = "'  -2.07%  "
$ws.Range("D30").Value = "'1.74"
$ws.Range("E30").Value = "'  -5.15%  "
$ws.Range("E31").Value = "'  +0.06%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.52"
$ws.Range("E32").Value = "'  -4.53%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'1.98"
$ws.Range("E33").Value = "'  -6.71%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'26.82"
$ws.Range("E34").Value = "'  -1.00%  "
$ws.Range("D35").Value = "'3.650.63"
$ws.Range("E35").Value = "'  +0.00%  "
$ws.Range("D36").Value = "'0.157"
$ws.Range("E36").Value = "'  -4.47%  "
$ws.Range("D37").Value = "'8.14"
$ws.Range("E37").Value = "'  -5.15%  "
$ws.Range("D38").Value = "'6.06"
$ws.Range("E38").Value = "'  -0.83%  "
$ws.Range("E39").Value = "'  +0.00%  "
$ws.Range("D40").Value = "'2.20"
$ws.Range("E40").Value = "'  +1.35%  "
$ws.Range("D41").Value = "'0.0896"
$ws.Range("E41").Value = "'  -5.34%  "
$ws.Range("E42").Value = "'  -0.06%  "
$ws.Range("E43").Value = "'  -2.13%  "
$ws.Range("D44").Value = "'165.86"
$ws.Range("E44").Value = "'  +4.84%  "
$ws.Range("D45").Value = "'47.52"
$ws.Range("E45").Value = "'  -1.36%  "
$ws.Range("D46").Value = "'0.000278"
$ws.Range("E46").Value = "'  -1.19%  "
$ws.Range("D47").Value = "'2.69"
$ws.Range("E47").Value = "'  -5.52%  "
$ws.Range("E48").Value = "'  +1.46%  "
$ws.Range("E49").Value = "'  -2.93%  "
$ws.Range("E50").Value = "'  -4.51%  "
$ws.Range("D51").Value = "'26.76"
$ws.Range("E51").Value = "'  -3.69%  "
